$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("更新日志")

# New log rows appended to the change log sheet (serial date numbers are
# used instead of Get-Date so the stored value is a clean integer with no
# time-of-day fraction). Row 60/61/62 each get their A:C values filled
# first; column D text is written afterwards, in the order row 60, then
# row 62, then row 61, so the rebuilt shared-string table gets the three
# new strings appended in the same order the source workbook has them
# (…, "删除用户功能", "CommonAction.class.php CommonModel.class.php",
# "附件上传。显示修改为Widget").
$rows = @(
    @{ Row = 60; Serial = 41788; B = "用户"; C = "新增"; D = "删除用户功能" },
    @{ Row = 61; Serial = 41791; B = "系统"; C = "改进"; D = "附件上传。显示修改为Widget" },
    @{ Row = 62; Serial = 41792; B = "系统"; C = "改进"; D = "CommonAction.class.php CommonModel.class.php" }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row
    $ws.Cells.Item(59, 1).Copy()
    $ws.Cells.Item($rowIndex, 1).PasteSpecial(-4122)
    $ws.Cells.Item($rowIndex, 1).Value = $r.Serial
    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Rows.Item($rowIndex).RowHeight = 21
}
$excel.CutCopyMode = $false

$ws.Cells.Item(60, 4).Value = "删除用户功能"
$ws.Cells.Item(62, 4).Value = "CommonAction.class.php CommonModel.class.php"
$ws.Cells.Item(61, 4).Value = "附件上传。显示修改为Widget"

$ws.Range("D67").Select() | Out-Null
